{"js": "// New values for the title paragraph and each table cell, in document order\n// (row-major: 20 rows x 5 columns).\nconst titleNew = \"2023-03-19 Sunday\";\nconst cellValues = [\n  [\"95-8=\", \"37+35=\", \"41+27=\", \"9+45=\", \"69-39=\"],\n  [\"25+71=\", \"80-47=\", \"58-35=\", \"67-24=\", \"61-59=\"],\n  [\"72-36=\", \"64+13=\", \"21-6=\", \"16-3=\", \"77-18=\"],\n  [\"12+34=\", \"83-34=\", \"78-47=\", \"5+26=\", \"77-34=\"],\n  [\"62-1=\", \"64+23=\", \"63-19=\", \"36+30=\", \"46+25=\"],\n  [\"18+3=\", \"80-32=\", \"20+20=\", \"44-5=\", \"94-21=\"],\n  [\"86+4=\", \"8+91=\", \"20-10=\", \"45-18=\", \"47-38=\"],\n  [\"57-35=\", \"21-8=\", \"21+17=\", \"0+23=\", \"97-49=\"],\n  [\"65-47=\", \"61-58=\", \"9+32=\", \"4+63=\", \"42-18=\"],\n  [\"52-29=\", \"76-40=\", \"58+41=\", \"16+13=\", \"58+4=\"],\n  [\"96-62=\", \"94-61=\", \"77-30=\", \"59-58=\", \"46+12=\"],\n  [\"50+5=\", \"78-9=\", \"27+67=\", \"13+65=\", \"17-2=\"],\n  [\"25+51=\", \"1-0=\", \"29+65=\", \"64-27=\", \"89-32=\"],\n  [\"73+5=\", \"9+21=\", \"25+46=\", \"76-5=\", \"48-30=\"],\n  [\"88-69=\", \"51+7=\", \"37+58=\", \"89-28=\", \"52-45=\"],\n  [\"75-62=\", \"67-47=\", \"56+21=\", \"37-12=\", \"63-11=\"],\n  [\"94-10=\", \"81-69=\", \"26-7=\", \"93-57=\", \"23+14=\"],\n  [\"91-0=\", \"68-50=\", \"74-72=\", \"42+4=\", \"74+4=\"],\n  [\"22+50=\", \"41-38=\", \"38+27=\", \"12-9=\", \"38+21=\"],\n  [\"97-27=\", \"19+35=\", \"33+29=\", \"21-17=\", \"40-27=\"],\n];\n\nconst body = context.document.body;\n\n// Update the date/title paragraph (first paragraph, outside the table).\nconst titlePara = body.paragraphs.getFirst();\ntitlePara.insertText(titleNew, \"Replace\");\n\n// Update each cell of the first (only) table, row by row, left to right.\nconst table = body.tables.getFirst();\nfor (let r = 0; r < cellValues.length; r++) {\n  for (let c = 0; c < cellValues[r].length; c++) {\n    table.getCell(r, c).value = cellValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# New values for the title paragraph and each table cell, in document order.\n$titleNew = \"2023-03-19 Sunday\"\n\n$cellValues = @(\n    \"95-8=\", \"37+35=\", \"41+27=\", \"9+45=\", \"69-39=\",\n    \"25+71=\", \"80-47=\", \"58-35=\", \"67-24=\", \"61-59=\",\n    \"72-36=\", \"64+13=\", \"21-6=\", \"16-3=\", \"77-18=\",\n    \"12+34=\", \"83-34=\", \"78-47=\", \"5+26=\", \"77-34=\",\n    \"62-1=\", \"64+23=\", \"63-19=\", \"36+30=\", \"46+25=\",\n    \"18+3=\", \"80-32=\", \"20+20=\", \"44-5=\", \"94-21=\",\n    \"86+4=\", \"8+91=\", \"20-10=\", \"45-18=\", \"47-38=\",\n    \"57-35=\", \"21-8=\", \"21+17=\", \"0+23=\", \"97-49=\",\n    \"65-47=\", \"61-58=\", \"9+32=\", \"4+63=\", \"42-18=\",\n    \"52-29=\", \"76-40=\", \"58+41=\", \"16+13=\", \"58+4=\",\n    \"96-62=\", \"94-61=\", \"77-30=\", \"59-58=\", \"46+12=\",\n    \"50+5=\", \"78-9=\", \"27+67=\", \"13+65=\", \"17-2=\",\n    \"25+51=\", \"1-0=\", \"29+65=\", \"64-27=\", \"89-32=\",\n    \"73+5=\", \"9+21=\", \"25+46=\", \"76-5=\", \"48-30=\",\n    \"88-69=\", \"51+7=\", \"37+58=\", \"89-28=\", \"52-45=\",\n    \"75-62=\", \"67-47=\", \"56+21=\", \"37-12=\", \"63-11=\",\n    \"94-10=\", \"81-69=\", \"26-7=\", \"93-57=\", \"23+14=\",\n    \"91-0=\", \"68-50=\", \"74-72=\", \"42+4=\", \"74+4=\",\n    \"22+50=\", \"41-38=\", \"38+27=\", \"12-9=\", \"38+21=\",\n    \"97-27=\", \"19+35=\", \"33+29=\", \"21-17=\", \"40-27=\"\n)\n\n# Update the date/title paragraph (first paragraph, outside the table).\n$d.Paragraphs.Item(1).Range.Text = $titleNew\n\n# Update each cell of the first (only) table, row by row, left to right.\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $cellValues[$i]\n        $i++\n    }\n}\n"}
